$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mdk"
$ws.Range("C2").Value = "Tspan1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 2.512729
$ws.Range("H2").Value = 7.538187
$ws.Range("I2").Value = 0.02190726325199687
$ws.Range("J2").Value = 0.02190726325199687
$ws.Range("K2").Value = 1.0
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.1177623333333333
$ws.Range("N2").Value = 0.353287
$ws.Range("O2").Value = 0.05622336632327191
$ws.Range("P2").Value = 0.05622336632327191
$ws.Range("Q2").Value = 0.2959048300743333
$ws.Range("R2").Value = 2.663143470669
$ws.Range("S2").Value = 0.001231700086957373
$ws.Range("T2").Value = 0.001231700086957373

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mdk"
$ws.Range("C3").Value = "Tspan1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 2.512729
$ws.Range("H3").Value = 7.538187
$ws.Range("I3").Value = 0.02190726325199687
$ws.Range("J3").Value = 0.02190726325199687
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 0.9102196666666668
$ws.Range("N3").Value = 2.730659
$ws.Range("O3").Value = 0.4345669137583306
$ws.Range("P3").Value = 0.4345669137583306
$ws.Range("Q3").Value = 2.287135352803667
$ws.Range("R3").Value = 20.584218175233
$ws.Range("S3").Value = 0.009520171780311568
$ws.Range("T3").Value = 0.00952017178031157

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Mdk"
$ws.Range("C4").Value = "Tspan1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 2.512729
$ws.Range("H4").Value = 7.538187
$ws.Range("I4").Value = 0.02190726325199687
$ws.Range("J4").Value = 0.02190726325199687
$ws.Range("K4").Value = 2.0
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1653493333333333
$ws.Range("N4").Value = 0.496048
$ws.Range("O4").Value = 0.07894286633226351
$ws.Range("P4").Value = 0.07894286633226351
$ws.Range("Q4").Value = 0.4154780649973333
$ws.Range("R4").Value = 3.739302584976
$ws.Range("S4").Value = 0.001729422154608097
$ws.Range("T4").Value = 0.001729422154608097

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Mdk"
$ws.Range("C5").Value = "Tspan1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 2.512729
$ws.Range("H5").Value = 7.538187
$ws.Range("I5").Value = 0.02190726325199687
$ws.Range("J5").Value = 0.02190726325199687
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 0.9012129999999999
$ws.Range("N5").Value = 2.703639
$ws.Range("O5").Value = 0.430266853586134
$ws.Range("P5").Value = 0.430266853586134
$ws.Range("Q5").Value = 2.264504040277
$ws.Range("R5").Value = 20.380536362493
$ws.Range("S5").Value = 0.009425969230119829
$ws.Range("T5").Value = 0.00942596923011983

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Mdk"
$ws.Range("C6").Value = "Tspan1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 92.89399466666667
$ws.Range("H6").Value = 278.681984
$ws.Range("I6").Value = 0.8098976036382196
$ws.Range("J6").Value = 0.8098976036382197
$ws.Range("K6").Value = 1.0
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.1177623333333333
$ws.Range("N6").Value = 0.353287
$ws.Range("O6").Value = 0.05622336632327191
$ws.Range("P6").Value = 0.05622336632327191
$ws.Range("Q6").Value = 10.93941356460089
$ws.Range("R6").Value = 98.45472208140801
$ws.Range("S6").Value = 0.04553516965369169
$ws.Range("T6").Value = 0.0455351696536917

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Mdk"
$ws.Range("C7").Value = "Tspan1"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 92.89399466666667
$ws.Range("H7").Value = 278.681984
$ws.Range("I7").Value = 0.8098976036382196
$ws.Range("J7").Value = 0.8098976036382197
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 0.9102196666666668
$ws.Range("N7").Value = 2.730659
$ws.Range("O7").Value = 0.4345669137583306
$ws.Range("P7").Value = 0.4345669137583306
$ws.Range("Q7").Value = 84.55394086082846
$ws.Range("R7").Value = 760.985467747456
$ws.Range("S7").Value = 0.3519547020733288
$ws.Range("T7").Value = 0.3519547020733289

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Mdk"
$ws.Range("C8").Value = "Tspan1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 92.89399466666667
$ws.Range("H8").Value = 278.681984
$ws.Range("I8").Value = 0.8098976036382196
$ws.Range("J8").Value = 0.8098976036382197
$ws.Range("K8").Value = 2.0
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1653493333333333
$ws.Range("N8").Value = 0.496048
$ws.Range("O8").Value = 0.07894286633226351
$ws.Range("P8").Value = 0.07894286633226351
$ws.Range("Q8").Value = 15.35996008880356
$ws.Range("R8").Value = 138.239640799232
$ws.Range("S8").Value = 0.0639356382668325
$ws.Range("T8").Value = 0.06393563826683252

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Mdk"
$ws.Range("C9").Value = "Tspan1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 92.89399466666667
$ws.Range("H9").Value = 278.681984
$ws.Range("I9").Value = 0.8098976036382196
$ws.Range("J9").Value = 0.8098976036382197
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 0.9012129999999999
$ws.Range("N9").Value = 2.703639
$ws.Range("O9").Value = 0.430266853586134
$ws.Range("P9").Value = 0.430266853586134
$ws.Range("Q9").Value = 83.71727561553067
$ws.Range("R9").Value = 753.455480539776
$ws.Range("S9").Value = 0.3484720936443666
$ws.Range("T9").Value = 0.3484720936443667

$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Mdk"
$ws.Range("C10").Value = "Tspan1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 1.610639333333333
$ws.Range("H10").Value = 4.831918
$ws.Range("I10").Value = 0.0140423817607685
$ws.Range("J10").Value = 0.0140423817607685
$ws.Range("K10").Value = 1.0
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.1177623333333333
$ws.Range("N10").Value = 0.353287
$ws.Range("O10").Value = 0.05622336632327191
$ws.Range("P10").Value = 0.05622336632327191
$ws.Range("Q10").Value = 0.1896726460517778
$ws.Range("R10").Value = 1.707053814466
$ws.Range("S10").Value = 0.0007895099737869194
$ws.Range("T10").Value = 0.0007895099737869195

$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Mdk"
$ws.Range("C11").Value = "Tspan1"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 1.610639333333333
$ws.Range("H11").Value = 4.831918
$ws.Range("I11").Value = 0.0140423817607685
$ws.Range("J11").Value = 0.0140423817607685
$ws.Range("K11").Value = 3.0
$ws.Range("L11").Value = 1.0
$ws.Range("M11").Value = 0.9102196666666668
$ws.Range("N11").Value = 2.730659
$ws.Range("O11").Value = 0.4345669137583306
$ws.Range("P11").Value = 0.4345669137583306
$ws.Range("Q11").Value = 1.466035597106889
$ws.Range("R11").Value = 13.194320373962
$ws.Range("S11").Value = 0.006102354503593439
$ws.Range("T11").Value = 0.00610235450359344

$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Mdk"
$ws.Range("C12").Value = "Tspan1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 1.610639333333333
$ws.Range("H12").Value = 4.831918
$ws.Range("I12").Value = 0.0140423817607685
$ws.Range("J12").Value = 0.0140423817607685
$ws.Range("K12").Value = 2.0
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.1653493333333333
$ws.Range("N12").Value = 0.496048
$ws.Range("O12").Value = 0.07894286633226351
$ws.Range("P12").Value = 0.07894286633226351
$ws.Range("Q12").Value = 0.2663181400071111
$ws.Range("R12").Value = 2.396863260064
$ws.Range("S12").Value = 0.001108545866326963
$ws.Range("T12").Value = 0.001108545866326963

$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Mdk"
$ws.Range("C13").Value = "Tspan1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 1.610639333333333
$ws.Range("H13").Value = 4.831918
$ws.Range("I13").Value = 0.0140423817607685
$ws.Range("J13").Value = 0.0140423817607685
$ws.Range("K13").Value = 3.0
$ws.Range("L13").Value = 1.0
$ws.Range("M13").Value = 0.9012129999999999
$ws.Range("N13").Value = 2.703639
$ws.Range("O13").Value = 0.430266853586134
$ws.Range("P13").Value = 0.430266853586134
$ws.Range("Q13").Value = 1.451529105511333
$ws.Range("R13").Value = 13.063761949602
$ws.Range("S13").Value = 0.006041971417061178
$ws.Range("T13").Value = 0.006041971417061179

$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Mdk"
$ws.Range("C14").Value = "Tspan1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 17.68108066666667
$ws.Range("H14").Value = 53.04324200000001
$ws.Range("I14").Value = 0.154152751349015
$ws.Range("J14").Value = 0.154152751349015
$ws.Range("K14").Value = 1.0
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.1177623333333333
$ws.Range("N14").Value = 0.353287
$ws.Range("O14").Value = 0.05622336632327191
$ws.Range("P14").Value = 0.05622336632327191
$ws.Range("Q14").Value = 2.082165315161556
$ws.Range("R14").Value = 18.739487836454
$ws.Range("S14").Value = 0.008666986608835917
$ws.Range("T14").Value = 0.008666986608835917

$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Mdk"
$ws.Range("C15").Value = "Tspan1"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 17.68108066666667
$ws.Range("H15").Value = 53.04324200000001
$ws.Range("I15").Value = 0.154152751349015
$ws.Range("J15").Value = 0.154152751349015
$ws.Range("K15").Value = 3.0
$ws.Range("L15").Value = 1.0
$ws.Range("M15").Value = 0.9102196666666668
$ws.Range("N15").Value = 2.730659
$ws.Range("O15").Value = 0.4345669137583306
$ws.Range("P15").Value = 0.4345669137583306
$ws.Range("Q15").Value = 16.09366735071978
$ws.Range("R15").Value = 144.843006156478
$ws.Range("S15").Value = 0.06698968540109677
$ws.Range("T15").Value = 0.06698968540109677

$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Mdk"
$ws.Range("C16").Value = "Tspan1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 17.68108066666667
$ws.Range("H16").Value = 53.04324200000001
$ws.Range("I16").Value = 0.154152751349015
$ws.Range("J16").Value = 0.154152751349015
$ws.Range("K16").Value = 2.0
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1653493333333333
$ws.Range("N16").Value = 0.496048
$ws.Range("O16").Value = 0.07894286633226351
$ws.Range("P16").Value = 0.07894286633226351
$ws.Range("Q16").Value = 2.923554900846223
$ws.Range("R16").Value = 26.311994107616
$ws.Range("S16").Value = 0.01216926004449594
$ws.Range("T16").Value = 0.01216926004449594

$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Mdk"
$ws.Range("C17").Value = "Tspan1"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3.0
$ws.Range("F17").Value = 1.0
$ws.Range("G17").Value = 17.68108066666667
$ws.Range("H17").Value = 53.04324200000001
$ws.Range("I17").Value = 0.154152751349015
$ws.Range("J17").Value = 0.154152751349015
$ws.Range("K17").Value = 3.0
$ws.Range("L17").Value = 1.0
$ws.Range("M17").Value = 0.9012129999999999
$ws.Range("N17").Value = 2.703639
$ws.Range("O17").Value = 0.430266853586134
$ws.Range("P17").Value = 0.430266853586134
$ws.Range("Q17").Value = 15.93441975084867
$ws.Range("R17").Value = 143.409777757638
$ws.Range("S17").Value = 0.06632681929458635
$ws.Range("T17").Value = 0.06632681929458635

